{"js": "// Replace the date line and each \"NNN\u00f7N=\" problem with its updated value.\n// Each (before, after) pair is applied as a single, case-sensitive,\n// first-match search & replace, in document order, so the identical\n// \"112\u00f72=\" text that appears twice across the whole edit (once removed,\n// once introduced) never collides with itself.\nconst replacements = [\n  [\"2025-09-12 Friday\", \"2025-09-13 Saturday\"],\n  [\"112\u00f72=\", \"113\u00f77=\"],\n  [\"649\u00f78=\", \"562\u00f78=\"],\n  [\"747\u00f76=\", \"726\u00f78=\"],\n  [\"629\u00f79=\", \"663\u00f74=\"],\n  [\"886\u00f74=\", \"755\u00f79=\"],\n  [\"217\u00f79=\", \"910\u00f75=\"],\n  [\"510\u00f77=\", \"885\u00f74=\"],\n  [\"785\u00f72=\", \"962\u00f72=\"],\n  [\"673\u00f73=\", \"492\u00f74=\"],\n  [\"769\u00f72=\", \"183\u00f73=\"],\n  [\"529\u00f74=\", \"112\u00f72=\"],\n  [\"531\u00f73=\", \"615\u00f79=\"],\n  [\"165\u00f79=\", \"542\u00f72=\"],\n  [\"440\u00f76=\", \"617\u00f73=\"],\n  [\"400\u00f79=\", \"968\u00f78=\"],\n  [\"743\u00f78=\", \"333\u00f75=\"],\n  [\"793\u00f76=\", \"124\u00f76=\"],\n  [\"311\u00f72=\", \"655\u00f73=\"],\n  [\"621\u00f73=\", \"603\u00f74=\"],\n  [\"540\u00f79=\", \"158\u00f76=\"],\n  [\"453\u00f76=\", \"111\u00f76=\"],\n  [\"393\u00f78=\", \"978\u00f76=\"],\n  [\"477\u00f76=\", \"124\u00f74=\"],\n  [\"317\u00f75=\", \"580\u00f72=\"],\n  [\"227\u00f75=\", \"634\u00f77=\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${before}`);\n  }\n  results.items[0].insertText(after, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"NNN\u00f7N=\" problem with its updated value.\n# Each (before, after) pair is unique in the document, so a scoped\n# Find/Replace (wdReplaceOne, first hit) for each pair applied in order is\n# safe even though some \"after\" values reintroduce a string that existed\n# (and was already replaced) earlier in the document.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-09-12 Friday\", \"2025-09-13 Saturday\"),\n    @(\"112\u00f72=\", \"113\u00f77=\"),\n    @(\"649\u00f78=\", \"562\u00f78=\"),\n    @(\"747\u00f76=\", \"726\u00f78=\"),\n    @(\"629\u00f79=\", \"663\u00f74=\"),\n    @(\"886\u00f74=\", \"755\u00f79=\"),\n    @(\"217\u00f79=\", \"910\u00f75=\"),\n    @(\"510\u00f77=\", \"885\u00f74=\"),\n    @(\"785\u00f72=\", \"962\u00f72=\"),\n    @(\"673\u00f73=\", \"492\u00f74=\"),\n    @(\"769\u00f72=\", \"183\u00f73=\"),\n    @(\"529\u00f74=\", \"112\u00f72=\"),\n    @(\"531\u00f73=\", \"615\u00f79=\"),\n    @(\"165\u00f79=\", \"542\u00f72=\"),\n    @(\"440\u00f76=\", \"617\u00f73=\"),\n    @(\"400\u00f79=\", \"968\u00f78=\"),\n    @(\"743\u00f78=\", \"333\u00f75=\"),\n    @(\"793\u00f76=\", \"124\u00f76=\"),\n    @(\"311\u00f72=\", \"655\u00f73=\"),\n    @(\"621\u00f73=\", \"603\u00f74=\"),\n    @(\"540\u00f79=\", \"158\u00f76=\"),\n    @(\"453\u00f76=\", \"111\u00f76=\"),\n    @(\"393\u00f78=\", \"978\u00f76=\"),\n    @(\"477\u00f76=\", \"124\u00f74=\"),\n    @(\"317\u00f75=\", \"580\u00f72=\"),\n    @(\"227\u00f75=\", \"634\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    # MatchCase:=True, MatchWholeWord:=True, Forward:=True, Wrap:=wdFindContinue(1),\n    # Replace:=wdReplaceOne(1) \u2014 replace just this (unique) occurrence.\n    $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 1)\n}\n"}
